# feat: add 2022-Q1 data
#
# This script:
#  1) Deletes the old "总计" (Total) sheet (its data is re-created afterwards)
#     so that the subsequent new-sheet sheetId allocation starts fresh.
#  2) Inserts a new worksheet "2022-Q1" right after "2021-Q4" (copying the
#     cell formatting from "2021-Q4") and fills it with the fund data.
#  3) Inserts a (new) "总计" worksheet right after "2022-Q1" (copying cell
#     formatting from "2021-Q4") and fills it with the updated summary,
#     which now includes a "2022-Q1" row on top of the previously existing
#     rows (2021-Q4 .. 2021-Q1).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 1: remove the existing "总计" sheet. We captured its (static,
# unchanged) data above in the plan and will re-write it below. Deleting
# it first means the sheetId counter allocates 5/6 to the two new sheets,
# matching the target sheet order/ids.
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

# ---------------------------------------------------------------------
# Step 2: create the new "2022-Q1" sheet right after "2021-Q4".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newQ1 = $wb.Worksheets.Add($null, $q4)
$newQ1.Name = "2022-Q1"

# Re-fetch references by name (structural changes can shift old handles).
$q4 = $wb.Worksheets.Item("2021-Q4")
$newQ1 = $wb.Worksheets.Item("2022-Q1")

# Copy cell formatting (styles/borders/fonts) from "2021-Q4" so the new
# sheet matches the look of the other quarterly sheets. Column A has no
# header cell on these sheets, so copy the header (B1:H1) and body
# (A2:H3) format ranges separately to avoid materializing a spurious A1
# cell.
$q4.Range("B1:H1").Copy()
$newQ1.Range("B1").PasteSpecial($xlPasteFormats)
$q4.Range("A2:H3").Copy()
$newQ1.Range("A2").PasteSpecial($xlPasteFormats)

# --- Header row ---
$newQ1.Range("B1").Value = "基金代码"
$newQ1.Range("C1").Value = "基金名称"
$newQ1.Range("D1").Value = "基金规模"
$newQ1.Range("E1").Value = "股票总仓位"
$newQ1.Range("F1").Value = "仓位占比"
$newQ1.Range("G1").Value = "持有市值(亿元)"
$newQ1.Range("H1").Value = "仓位排名"

# --- Row 2: 010343 / 华宝英国富时100指数（QDII）A ---
$newQ1.Range("A2").Value = 0
$newQ1.Range("B2").NumberFormat = "@"
$newQ1.Range("B2").Value = "010343"
$newQ1.Range("B2").Style = "Normal"
$newQ1.Range("C2").Value = "华宝英国富时100指数（QDII）A"
$newQ1.Range("D2").NumberFormat = "@"
$newQ1.Range("D2").Value = "0.20"
$newQ1.Range("D2").Style = "Normal"
$newQ1.Range("E2").NumberFormat = "@"
$newQ1.Range("E2").Value = "93.65"
$newQ1.Range("E2").Style = "Normal"
$newQ1.Range("F2").NumberFormat = "@"
$newQ1.Range("F2").Value = "3.85"
$newQ1.Range("F2").Style = "Normal"
$newQ1.Range("G2").NumberFormat = "@"
$newQ1.Range("G2").Value = "0.0077"
$newQ1.Range("G2").Style = "Normal"
$newQ1.Range("H2").Value = 6

# --- Row 3: 010344 / 华宝英国富时100指数（QDII）C ---
$newQ1.Range("A3").Value = 1
$newQ1.Range("B3").NumberFormat = "@"
$newQ1.Range("B3").Value = "010344"
$newQ1.Range("B3").Style = "Normal"
$newQ1.Range("C3").Value = "华宝英国富时100指数（QDII）C"
$newQ1.Range("D3").NumberFormat = "@"
$newQ1.Range("D3").Value = "0.06"
$newQ1.Range("D3").Style = "Normal"
$newQ1.Range("E3").NumberFormat = "@"
$newQ1.Range("E3").Value = "93.65"
$newQ1.Range("E3").Style = "Normal"
$newQ1.Range("F3").NumberFormat = "@"
$newQ1.Range("F3").Value = "3.85"
$newQ1.Range("F3").Style = "Normal"
$newQ1.Range("G3").NumberFormat = "@"
$newQ1.Range("G3").Value = "0.0023"
$newQ1.Range("G3").Style = "Normal"
$newQ1.Range("H3").Value = 6

# ---------------------------------------------------------------------
# Step 3: create the new "总计" sheet right after "2022-Q1".
# ---------------------------------------------------------------------
$newQ1 = $wb.Worksheets.Item("2022-Q1")
$newTotal = $wb.Worksheets.Add($null, $newQ1)
$newTotal.Name = "总计"

# Re-fetch references by name.
$q4 = $wb.Worksheets.Item("2021-Q4")
$newTotal = $wb.Worksheets.Item("总计")

# Copy cell formatting from "2021-Q4": the header look (B1:D1) and the
# A-column style (single A2 cell, extended down through row 6). Column A
# has no header cell, so this avoids materializing a spurious A1 cell.
$q4.Range("B1:D1").Copy()
$newTotal.Range("B1").PasteSpecial($xlPasteFormats)
$q4.Range("A2").Copy()
$newTotal.Range("A2:A6").PasteSpecial($xlPasteFormats)

# --- Header row ---
$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

# --- Row 2 (new): 2022-Q1 ---
$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 2
$newTotal.Range("D2").Value = 0.01

# --- Row 3: 2021-Q4 ---
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 2
$newTotal.Range("D3").Value = 0.01

# --- Row 4: 2021-Q3 ---
$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 2
$newTotal.Range("D4").Value = 0.01

# --- Row 5: 2021-Q2 ---
$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 2
$newTotal.Range("D5").Value = 0.01

# --- Row 6: 2021-Q1 ---
$newTotal.Range("A6").Value = 4
$newTotal.Range("B6").Value = "2021-Q1"
$newTotal.Range("C6").Value = 2
$newTotal.Range("D6").Value = 0.01
